$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRow = 45
$newRow = 46

# Copy the formatting (style) from the prior row so the new row matches
$ws.Range("A$srcRow`:H$srcRow").Copy() | Out-Null
$ws.Range("A$newRow`:H$newRow").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Cells.Item($newRow, 1).Value = "2025-08-23 03:48:00 UTC"
$ws.Cells.Item($newRow, 2).Value = "2025-08-23 09:18:00 IST"
$ws.Cells.Item($newRow, 3).Value = "SKIPPED"
$ws.Cells.Item($newRow, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($newRow, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$ws.Cells.Item($newRow, 6).Value = ""
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = ""
